# B2 R3 bottle mass
# Fill in rain_date (C) and bottle_mass(g) (E) for rows 20-37 on the
# "Rainfall 3" worksheet, and move the active selection to E21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rainfall 3")

$rainDate = 20210616

$bottleMass = @{
    20 = 89
    21 = 89
    22 = 88
    23 = 89
    24 = 89
    25 = 88
    26 = 104
    27 = 88
    28 = 89
    29 = 87
    30 = 101
    31 = 90
    32 = 89
    33 = 88
    34 = 89
    35 = 88
    36 = 89
    37 = 89
}

foreach ($row in 20..37) {
    $ws.Cells.Item($row, 3).Value = $rainDate
    $ws.Cells.Item($row, 5).Value = $bottleMass[$row]
}

$ws.Range("E21").Select()
